$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its textual formatting so that
# values such as "1.001" or "18.80" are not reinterpreted as numbers
# and lose significant trailing zeros / thousands separators.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.016.12'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '1.901.55'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '0.7441'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '242.21'
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '0.3072'
$ws.Range("E8").Value = '  -2.42%  '
$ws.Range("D9").Value = '25.61'
$ws.Range("E9").Value = '  -6.81%  '
$ws.Range("D10").Value = '0.06897'
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("D11").Value = '0.08030'
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '0.7549'
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("D13").Value = '1.907.47'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("D14").Value = '5.235'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = '91.23'
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").Value = '6.193'
$ws.Range("E16").Value = '  +3.97%  '
$ws.Range("D17").Value = '30.029.46'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = '14.04'
$ws.Range("E18").Value = '  -2.93%  '
$ws.Range("D19").Value = '0.000007782'
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").Value = '237.11'
$ws.Range("E20").Value = '  -5.10%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = '2.152.55'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = '7.101'
$ws.Range("E24").Value = '  +6.81%  '
$ws.Range("D25").Value = '9.347'
$ws.Range("E25").Value = '  -2.15%  '
$ws.Range("D26").Value = '167.81'
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").Value = '18.80'
$ws.Range("E27").Value = '  -1.14%  '
$ws.Range("D28").Value = '0.1272'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").Value = '2.047'
$ws.Range("E29").Value = '  -5.36%  '
$ws.Range("D30").Value = '1.351'
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("D31").Value = '1.534'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '4.306'
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").Value = '4.053'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").Value = '0.05306'
$ws.Range("E34").Value = '  +1.81%  '
$ws.Range("D35").Value = '1.284'
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("D36").Value = '0.7387'
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("D37").Value = '2.727'
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = '0.01948'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").Value = '2.767'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").Value = '6.269'
$ws.Range("E40").Value = '  -3.87%  '
$ws.Range("D41").Value = '0.4459'
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").Value = '72.68'
$ws.Range("E42").Value = '  -5.12%  '
$ws.Range("D43").Value = '1.950'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").Value = '7.758'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("D46").Value = '0.8323'
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("D47").Value = '101.55'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '9.873'
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").Value = '2.055.20'
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = '36.62'
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("D51").Value = '0.05990'
$ws.Range("E51").Value = '  -0.42%  '
